# Mike coding while frustrated...
# Transpose the ticker header row (A1:V1, "Date" + 21 tickers) into a single
# "Symbol" column (A1:A22) running down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing ticker symbols (everything after the "Date" label in A1)
$lastCol = $ws.Cells.Item(1, 1).End(-4161).Column  # xlToRight = -4161
$tickers = @()
for ($c = 2; $c -le $lastCol; $c++) {
    $tickers += $ws.Cells.Item(1, $c).Value()
}

# Wipe out the old horizontal layout
$ws.Range("A1:V1").ClearContents()

# Write the new vertical layout: header "Symbol" in A1, tickers down column A
$ws.Cells.Item(1, 1).Value = "Symbol"
$row = 2
foreach ($ticker in $tickers) {
    $ws.Cells.Item($row, 1).Value = $ticker
    $row++
}

$ws.Range("C7").Select()
